# Program_List.xlsx edit: "added code 49 to 55"
#
# - Corrects the "Written" status of 5 existing rows (Leap year .. Voter
#   eligiblity) from "i" (in-progress) to "p" (done).
# - Appends 12 new programs to Table1 (rows 50-61), each with a serial
#   number, title, Status="p"/"i" and Written="p"/"i" like the rest of
#   the sheet.
# - Re-points the page setup to print at 46% scale without a fixed page
#   height, and leaves the sheet scrolled to/selecting the newly typed
#   area.

$wb  = $excel.ActiveWorkbook
$ws  = $wb.Worksheets.Item(1)
$tbl = $ws.ListObjects.Item(1)

# ---------------------------------------------------------------------------
# 1) Five already-present rows get their "Written" column corrected from
#    "i" to "p".
# ---------------------------------------------------------------------------
foreach ($r in 15..19) {
    $ws.Cells.Item($r, 4).Value = "p"
}

# ---------------------------------------------------------------------------
# 2) Append 12 new data rows to the table (sheet rows 50-61).
# ---------------------------------------------------------------------------
$newRows = @(
    @{ Title = "Copy constructor";                                                  Status = "p"; Written = "i" },
    @{ Title = "Decimal to binary conversion";                                      Status = "p"; Written = "i" },
    @{ Title = "Nelson Number (111,222,333)";                                       Status = "p"; Written = "i" },
    @{ Title = "Operater overloading usind friend function";                        Status = "p"; Written = "i" },
    @{ Title = "Find the day from date of birth";                                   Status = "p"; Written = "i" },
    @{ Title = "Matrix Row and column sum program";                                 Status = "p"; Written = "i" },
    @{ Title = "String Function -> user defined strlen,strcpy,strcat,strcmp";       Status = "p"; Written = "i" },
    @{ Title = "1D array ( INSERTION,DELETION,TRAVERSAL)";                          Status = "i"; Written = "i" },
    @{ Title = "QUEUE ( INSERTION,DELETION,TRAVERSAL)";                             Status = "i"; Written = "i" },
    @{ Title = "STACK ( INSERTION,DELETION,TRAVERSAL)";                             Status = "i"; Written = "i" },
    @{ Title = "LINKED LIST ( INSERTION,DELETION,TRAVERSAL)";                       Status = "i"; Written = "i" },
    @{ Title = "TREE ( INSERTION,DELETION,TRAVERSAL)";                              Status = "i"; Written = "i" }
)

$lastRow = $tbl.Range.Row + $tbl.Range.Rows.Count - 1
$serial  = $ws.Cells.Item($lastRow, 1).Value()

foreach ($row in $newRows) {
    $tbl.ListRows.Add() | Out-Null
    $lastRow = $tbl.Range.Row + $tbl.Range.Rows.Count - 1
    $serial  = $serial + 1

    # Copy formatting (styles + row height) from the row directly above so
    # the new row matches the rest of the table.
    $ws.Range("A" + ($lastRow - 1) + ":D" + ($lastRow - 1)).Copy() | Out-Null
    $ws.Range("A" + $lastRow + ":D" + $lastRow).PasteSpecial(-4122) | Out-Null
    $ws.Rows.Item($lastRow).RowHeight = $ws.Rows.Item($lastRow - 1).RowHeight

    $ws.Cells.Item($lastRow, 1).Value = $serial
    $ws.Cells.Item($lastRow, 2).Value = $row.Title
    $ws.Cells.Item($lastRow, 3).Value = $row.Status
    $ws.Cells.Item($lastRow, 4).Value = $row.Written
}

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 3) Extend the "i"/"p" highlighting conditional formatting so it also
#    covers the Status/Written cells of the new rows.
# ---------------------------------------------------------------------------
$firstNewRow = $lastRow - $newRows.Count + 1
$redFont   = 393372     # RGB(0x9C,0x00,0x06) in COM BGR order
$redFill   = 13551615   # RGB(0xFF,0xC7,0xCE)
$greenFont = 24832      # RGB(0x00,0x61,0x00)
$greenFill = 13561798   # RGB(0xC6,0xEF,0xCE)

$statusRange  = $ws.Range("C" + $firstNewRow + ":C" + $lastRow)
$writtenRange = $ws.Range("D" + $firstNewRow + ":D" + $lastRow)

foreach ($rng in @($statusRange, $writtenRange)) {
    $fcI = $rng.FormatConditions.Add(1, 3, '"i"')
    $fcI.Font.Color = $redFont
    $fcI.Interior.Color = $redFill

    $fcP = $rng.FormatConditions.Add(1, 3, '"p"')
    $fcP.Font.Color = $greenFont
    $fcP.Interior.Color = $greenFill
}

# ---------------------------------------------------------------------------
# 4) Page setup: scale printout to 46% and drop the fixed page-height
#    constraint; select the newly entered area.
# ---------------------------------------------------------------------------
$ws.PageSetup.Zoom = 46
$ws.PageSetup.FitToPagesTall = 0

$ws.Range("B57").Select() | Out-Null
